# repull data, push all data, mean calculation
# Update the dSF column (F) values for the affected rows to reflect the
# repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -4
    4  = 3
    5  = -10
    9  = -4
    13 = 5
    21 = 4
    22 = -2
    23 = -2
    24 = -5
    28 = -5
    32 = -4
    36 = -1
    39 = -2
    46 = -2
    56 = 3
    60 = -3
    61 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
